$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.535.30"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "2.928.31"
$ws.Range("E3").Value = "  -2.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "376.33"
$ws.Range("E5").Value = "  +6.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.23"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.96"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0837"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.37"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").Value = "3.398.39"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.37"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("D16").Value = "2.934.58"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.940"
$ws.Range("E17").Value = "  -7.26%  "
$ws.Range("D18").Value = "51.503.09"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.00"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.41"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.10"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.78"
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("E26").Value = "  -5.40%  "
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.82"
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.03"
$ws.Range("E31").Value = "  +8.66%  "
$ws.Range("E32").Value = "  -5.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.83"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "34.15"
$ws.Range("E36").Value = "  -4.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0427"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").Value = "  -9.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.98"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("E41").Value = "  -6.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -6.14%  "
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.46"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.86"
$ws.Range("E45").Value = "  -6.09%  "
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.275"
$ws.Range("E47").Value = "  +13.19%  "
$ws.Range("D48").Value = "2.023.90"
$ws.Range("E48").Value = "  -4.84%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.17"
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("D51").Value = "3.219.57"
$ws.Range("E51").Value = "  -2.61%  "
